$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions) - update column F ("想去人数" / interest count)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1190
$ws1.Range("F4").Value = 278
$ws1.Range("F6").Value = 14
$ws1.Range("F7").Value = 12335
$ws1.Range("F8").Value = 61
$ws1.Range("F9").Value = 17
$ws1.Range("F10").Value = 13
$ws1.Range("F11").Value = 159
$ws1.Range("F12").Value = 12147
$ws1.Range("F13").Value = 4832
$ws1.Range("F14").Value = 4695
$ws1.Range("F15").Value = 130
$ws1.Range("F16").Value = 64
$ws1.Range("F17").Value = 422
$ws1.Range("F19").Value = 950
$ws1.Range("F20").Value = 3

# Sheet "全部类型" (All types) - mirrors the same data, update column F
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1190
$ws4.Range("F4").Value = 278
$ws4.Range("F8").Value = 14
$ws4.Range("F9").Value = 12335
$ws4.Range("F10").Value = 61
$ws4.Range("F11").Value = 17
$ws4.Range("F12").Value = 13
$ws4.Range("F13").Value = 159
$ws4.Range("F14").Value = 12147
$ws4.Range("F15").Value = 4832
$ws4.Range("F16").Value = 4695
$ws4.Range("F17").Value = 130
$ws4.Range("F18").Value = 64
$ws4.Range("F19").Value = 422
$ws4.Range("F21").Value = 950
$ws4.Range("F22").Value = 3

